$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Förändrad" (column C) date from 45204 to 45205 for every existing
#    data row (2-500).
$ws.Range("C2:C500").Value = 45205

# Row 500 picks up an explicit row height once a new row is appended after it
# (matching every other data row in the sheet).
$ws.Rows.Item(500).RowHeight = 15

# 2. Copy the formatting (date formats on B/C, wrap-text on R) from row 500
#    into the new row 501 before filling in its values. Columns are copied
#    individually (skipping F, which has no data on row 500/501) to avoid
#    introducing stray empty cells.
$ws.Range("B500").Copy()
$ws.Range("B501").PasteSpecial(-4122)
$ws.Range("C500").Copy()
$ws.Range("C501").PasteSpecial(-4122)
$ws.Range("R500").Copy()
$ws.Range("R501").PasteSpecial(-4122)

# 3. Fill in the new record (row 501).
$ws.Range("A501").Value = "A 47891-2023"
$ws.Range("B501").Value = 45204
$ws.Range("C501").Value = 45205
$ws.Range("D501").Value = "JÖNKÖPINGS LÄN"
$ws.Range("E501").Value = "GISLAVED"
$ws.Range("G501").Value = 4.4
$ws.Range("H501").Value = 0
$ws.Range("I501").Value = 0
$ws.Range("J501").Value = 0
$ws.Range("K501").Value = 0
$ws.Range("L501").Value = 0
$ws.Range("M501").Value = 0
$ws.Range("N501").Value = 0
$ws.Range("O501").Value = 0
$ws.Range("P501").Value = 0
$ws.Range("Q501").Value = 0
